$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Length -gt 0) { $t = $t.Substring(0, $t.Length - 1) }
        if ($t -eq $text) {
            return $i
        }
    }
    throw "paragraph with text '$text' not found"
}

# ---------------------------------------------------------------------------
# 1) "Possibly a leader board using TKint<bookmark>er" -> "...using TKinter"
#    Word's Find/Replace naturally spans the split runs (and the _GoBack
#    bookmark sitting between them), replacing the whole match with a single
#    contiguous run and dropping the now-redundant bookmark from that spot.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("TKinter", $false, $false, $false, $false, $false, `
    $true, 1, $false, "TKinter", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove the "Work -" heading paragraph and the blank paragraph after it.
# ---------------------------------------------------------------------------
$workIdx = Get-ParaIndexByText $d "Work –"
$pStart = $d.Paragraphs.Item($workIdx)
$pAfterBlank = $d.Paragraphs.Item($workIdx + 2)
$d.Range($pStart.Range.Start, $pAfterBlank.Range.Start).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3) "Patryk -  " paragraph collapses to a single space (tabs + name gone).
# ---------------------------------------------------------------------------
$patrykIdx = Get-ParaIndexByText $d "`tPatryk –  `t"
$pPatryk = $d.Paragraphs.Item($patrykIdx)
$rngPatryk = $pPatryk.Range
$rngPatryk.End = $rngPatryk.End - 1
$rngPatryk.Text = " "

# ---------------------------------------------------------------------------
# 4) Drop the Jamie / Tom / Euan / Erik paragraphs (and their blank spacer
#    paragraphs) entirely.
# ---------------------------------------------------------------------------
$jamieIdx = Get-ParaIndexByText $d "`tJamie - "
$conorIdx = Get-ParaIndexByText $d "`tConor - "
$pJamie = $d.Paragraphs.Item($jamieIdx)
$pConor = $d.Paragraphs.Item($conorIdx)
$d.Range($pJamie.Range.Start, $pConor.Range.Start).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 5) The Conor paragraph itself stays but becomes empty, and the _GoBack
#    bookmark now lives inside it.
# ---------------------------------------------------------------------------
$conorIdx = Get-ParaIndexByText $d "`tConor - "
$pConor = $d.Paragraphs.Item($conorIdx)
$rngConor = $pConor.Range
$rngConor.End = $rngConor.End - 1
$rngConor.Text = ""
$d.Bookmarks.Add("_GoBack", $rngConor) | Out-Null
